$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(30)
$r = $p.Range
$insPt = $d.Range($r.Start, $r.Start)
$xml = @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p><w:proofErr w:type='spellStart'/><w:r><w:t>dynd</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>==0.7.3.dev1</w:t></w:r></w:p>
<w:p><w:proofErr w:type='spellStart'/><w:r><w:lastRenderedPageBreak/><w:t>nb-anacondacloud</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>==1.1.0</w:t></w:r></w:p>
<w:p><w:proofErr w:type='spellStart'/><w:r><w:t>nb-conda</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>==1.1.0</w:t></w:r></w:p>
<w:p></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$insPt.InsertXML($xml) | Out-Null

$emptyPara = $d.Paragraphs.Item(33)
$emptyPara.Range.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("dynd==0.7.3.dev1", $false, $false, $false, $false, $false, $true, 1, $false, "nb-conda-kernels==1.0.3", 2) | Out-Null
Write-Output "done"
